# The source data gained a new weekly record for "Ajo" (Chino / Primera)
# at Terminal La Palmera de La Serena, Coquimbo. It belongs right after the
# existing row 488, so every subsequent record shifts down by one row and
# the new record's values are written into the freshly-inserted row 489.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 489, pushing rows 489:510 down to 490:511.
$ws.Rows.Item(489).Insert()

$newRow = @{
    A = 8
    B = "Terminal La Palmera de La Serena"
    C = "Coquimbo"
    D = 45147
    E = 4
    F = 100112003
    G = "Ajo"
    H = "Chino"
    I = "Primera"
    J = 360
    K = 21000
    L = 22000
    M = 21500
    N = '$/caja 10 kilos'
    O = "China"
    P = 2150
    Q = 10
    R = "Hortaliza"
}

$ws.Range("A489").Value = $newRow.A
$ws.Range("B489").Value = $newRow.B
$ws.Range("C489").Value = $newRow.C
$ws.Range("D489").Value = $newRow.D
$ws.Range("E489").Value = $newRow.E
$ws.Range("F489").Value = $newRow.F
$ws.Range("G489").Value = $newRow.G
$ws.Range("H489").Value = $newRow.H
$ws.Range("I489").Value = $newRow.I
$ws.Range("J489").Value = $newRow.J
$ws.Range("K489").Value = $newRow.K
$ws.Range("L489").Value = $newRow.L
$ws.Range("M489").Value = $newRow.M
$ws.Range("N489").Value = $newRow.N
$ws.Range("O489").Value = $newRow.O
$ws.Range("P489").Value = $newRow.P
$ws.Range("Q489").Value = $newRow.Q
$ws.Range("R489").Value = $newRow.R
